$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update N6 assertion value (now expects a result instead of "None")
$ws.Range("N6").Value = "data.id: contains 7"

# Row 7 becomes the "/register" POST (BasicAuth, 400) test case,
# with updated, more specific credentials text.
$ws.Range("C7").Value = "/register"
$ws.Range("D7").Value = "POST"
$ws.Range("H7").Value = "None"
$ws.Range("J7").Value = "BasicAuth"
$ws.Range("K7").Value = "username: eve.holt@reqres.in, password: cityslicka"
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "400"
$ws.Range("N7").Value = "None"

# Row 8 becomes the "/users/3" GET (200) test case.
$ws.Range("C8").Value = "/users/3"
$ws.Range("D8").Value = "GET"
$ws.Range("H8").Value = "id=3"
$ws.Range("J8").Value = "None"
$ws.Range("K8").Value = "None"
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "200"
$ws.Range("N8").Value = "data.last_name: 'Wong'"

# The trailing "/users/2" GET row (row 9) is removed entirely.
$ws.Rows("9").Delete()
